# "turn off profanity filter, adjust gui size"
#
# The swear-word list in column A gets two new entries inserted in
# alphabetical position:
#   - "putang ina"  -> right after "putangina" (was row 21, now pushes
#                       putragis..ungas down by one)
#   - "tang ina"    -> right after "tangina"   (pushes tarantado..ungas
#                       down by one more)
#
# Net effect on column A (rows are 1-based, row 1 is the header
# "swear_words_list"):
#   A21 putang ina   (new)
#   A22 putragis      (was A21)
#   A23 taena         (was A22)
#   A24 tanga         (was A23)
#   A25 tangina       (was A24)
#   A26 tang ina      (new)
#   A27 tarantado     (was A25)
#   A28 ulol          (was A26)
#   A29 ulul          (was A27)
#   A30 ungas         (was A28)
# and one previously-blank trailing row is pushed in, growing the sheet
# from A1:A34 to A1:A35.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 21-30 already share the same style as the rest of the list (s="1"),
# so the reshuffled values can simply be written in place - no row/style
# shifting is required for this part.
$newTail = @(
    "putang ina",
    "putragis",
    "taena",
    "tanga",
    "tangina",
    "tang ina",
    "tarantado",
    "ulol",
    "ulul",
    "ungas"
)

$row = 21
foreach ($word in $newTail) {
    $ws.Cells.Item($row, 1).Value = $word
    $row = $row + 1
}

# The trailing blank rows grow by one (A29:A34 -> A31:A35), so insert a
# single blank row just below the list to shift the rest of the blank
# rows down and extend the used range to A1:A35.
$ws.Range("A31").Insert()

# Adjust the on-screen view to match: scroll so row 16 is near the top
# and leave the selection on C27.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("C27").Select()
